$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("A1").Value = "sasdsfads"
